$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# New TOTAL row (row 4): label in B4 (right-aligned, new style), and
# the PRESUPUESTO / VENTA / POR CUMPLIR / CUMPLIMIENTO totals in C4:F4,
# reusing the number-format styles already used by rows 2-3.
$ws.Range("B4").Value = "TOTAL"
$ws.Range("B4").HorizontalAlignment = -4152

$ws.Range("C2").Copy($ws.Range("C4"))
$ws.Range("C4").Value = 17500

$ws.Range("D2").Copy($ws.Range("D4"))
$ws.Range("D4").Value = 605.48

$ws.Range("E2").Copy($ws.Range("E4"))
$ws.Range("E4").Value = 16894.52

$ws.Range("F2").Copy($ws.Range("F4"))
$ws.Range("F4").Value = 0.03459885714285715
